# Weekly fruit/vegetable price update.
# Insert a new weekly price record as row 28 (pushing the existing
# rows 28-33 down to 29-34) and populate it with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 28:33 down to 29:34, making room for the new record at row 28.
$ws.Rows.Item(28).Insert()

# Populate the new row 28 with this week's Perejil (parsley) price data.
$ws.Cells.Item(28, 1).Value = 1
$ws.Cells.Item(28, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(28, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(28, 4).Value = 44917
$ws.Cells.Item(28, 5).Value = 15
$ws.Cells.Item(28, 6).Value = 100112044
$ws.Cells.Item(28, 7).Value = "Perejil"
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 300
$ws.Cells.Item(28, 11).Value = 2700
$ws.Cells.Item(28, 12).Value = 3000
$ws.Cells.Item(28, 13).Value = 2850
$ws.Cells.Item(28, 14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(28, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(28, 16).Value = 1425
$ws.Cells.Item(28, 17).Value = 2
$ws.Cells.Item(28, 18).Value = "Hortaliza"
